$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rId1 / sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1861
$ws1.Range("G2").Value = 75
$ws1.Range("F3").Value = 277
$ws1.Range("F7").Value = 165
$ws1.Range("F8").Value = 3675
$ws1.Range("F9").Value = 162
$ws1.Range("F14").Value = 641
$ws1.Range("F16").Value = 833
$ws1.Range("F21").Value = 80
$ws1.Range("F23").Value = 3042
$ws1.Range("F24").Value = 5428
$ws1.Range("F29").Value = 3150
$ws1.Range("F30").Value = 327
$ws1.Range("F31").Value = 2340
$ws1.Range("F35").Value = 159
$ws1.Range("F36").Value = 217
$ws1.Range("F38").Value = 75
$ws1.Range("F40").Value = 841
$ws1.Range("F42").Value = 21
$ws1.Range("F44").Value = 52
$ws1.Range("F45").Value = 517

# --- Sheet "全部类型" (rId4 / sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1861
$ws4.Range("G2").Value = 75
$ws4.Range("F3").Value = 277
$ws4.Range("F7").Value = 165
$ws4.Range("F8").Value = 3675
$ws4.Range("F9").Value = 162
$ws4.Range("F15").Value = 641
$ws4.Range("F17").Value = 833
$ws4.Range("F22").Value = 80
$ws4.Range("F24").Value = 3042
$ws4.Range("F25").Value = 5428
$ws4.Range("F30").Value = 3150
$ws4.Range("F31").Value = 327
$ws4.Range("F32").Value = 2340
$ws4.Range("F36").Value = 159
$ws4.Range("F37").Value = 217
$ws4.Range("F39").Value = 75
$ws4.Range("F41").Value = 841
$ws4.Range("F43").Value = 21
$ws4.Range("F45").Value = 52
$ws4.Range("F46").Value = 517
